# The commit swaps the two theme parts in this deck:
#   ppt/theme/theme1.xml  (the slide master's theme -- "Integral" / Red
#                           Violet colour scheme) becomes the default
#                           Office "Office Theme" colour scheme.
#   ppt/theme/theme2.xml  (the notes master's theme -- "Office Theme")
#                           becomes the "Integral" / Red Violet colour
#                           scheme that used to live in theme1.xml.
#
# Font scheme ("Office") and format scheme ("Office") are identical
# between the two theme parts already, so only the 12 colour-scheme
# entries (and, where reachable, the theme's display name) actually
# need to change.
#
# The only theme object this COM host exposes is the one bound to the
# presentation's slide master (ppt/theme/theme1.xml), reached through
# Designs(1).SlideMaster.Theme / ActivePresentation.SlideMaster.Theme.
# We drive every one of its twelve ThemeColorScheme entries to the
# standard Office palette, in clrScheme document order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

$colors.Item(1).RGB  = 0         # dk1      000000
$colors.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388   # dk2      44546A
$colors.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501   # accent2  ED7D31
$colors.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Item(8).RGB  = 49407     # accent4  FFC000
$colors.Item(9).RGB  = 12874308  # accent5  4472C4
$colors.Item(10).RGB = 4697456   # accent6  70AD47
$colors.Item(11).RGB = 12673797  # hlink    0563C1
$colors.Item(12).RGB = 7491477   # folHlink 954F72
